# Uganda central resources - fix missing entity names
#
# The "Data" worksheet lists one row per district (identified by a code such
# as "d101" in column A) but column B ("entity-name") was left blank for every
# row. This script fills in column B with the corresponding Uganda district
# name for each of the 112 data rows (rows 2-113).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Data")

# District names, in the same row order as the existing id codes in column A
# (row 2 = d101 ... row 113 = d426).
$names = @(
    "Kalangala","Kampala","Kiboga","Luwero","Masaka","Mpigi","Mubende","Mukono",
    "Nakasongola","Rakai","Ssembabule","Kayunga","Wakiso","Lyantonde","Mityana",
    "Nakaseke","Buikwe","Bukomansimbi","Butambala","Buvuma","Gomba","Kalungu",
    "Kyankwanzi","Lwengo","Bugiri","Busia","Iganga","Jinja","Kamuli","Kapchorwa",
    "Katakwi","Kumi","Mbale","Pallisa","Soroti","Tororo","Kaberamaido","Mayuge",
    "Sironko","Amuria","Budaka","Bududa","Bukedea","Bukwo","Butaleja","Kaliro",
    "Manafwa","Namutumba","Bulambuli","Buyende","Kibuku","Kween","Luuka",
    "Namayingo","Ngora","Serere","Adjumani","Apac","Arua","Gulu","Kitgum",
    "Kotido","Lira","Moroto","Moyo","Nebbi","Nakapiripirit","Pader","Yumbe",
    "Abim","Amolatar","Amuru","Dokolo","Kaabong","Koboko","Maracha","Oyam",
    "Agago","Alebtong","Amudat","Kole","Lamwo","Napak","Nwoya","Otuke","Zombo",
    "Bundibugyo","Bushenyi","Hoima","Kabale","Kabarole","Kasese","Kibaale",
    "Kisoro","Masindi","Mbarara","Ntungamo","Rukungiri","Kamwenge","Kanungu",
    "Kyenjojo","Buliisa","Ibanda","Isingiro","Kiruhura","Buhweju","Kiryandongo",
    "Kyegegwa","Mitooma","Ntoroko","Rubirizi","Sheema"
)

$startRow = 2
for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $names[$i]
}
